# Update table summarizing transportation/manipulation approaches
# (refines Trajectory Planning / Agent Control Strategy / Measurements text,
#  plus a few Task/Load-Config/Centralized-Distributed corrections)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 'Calculation of workspace of valid trajectories, no info on how a trajectory is selected'
$ws.Range("G2").Value = 'PID + feedforward compensation + potential field (collision avoidance)'
$ws.Range("H2").Value = 'Agents'' state'
$ws.Range("F3").Value = 'Payload trajectory set at beginning'
$ws.Range("G3").Value = 'PD nonlinear'
$ws.Range("H3").Value = 'Agents'' state'
$ws.Range("F4").Value = 'Payload trajectory set at beginning (piecewise polynomial over time)'
$ws.Range("G4").Value = 'LQR (leader) +  LQR & PD (follower)'
$ws.Range("H4").Value = 'Agents'' state'
$ws.Range("F5").Value = 'Payload trajectory set at beginning (piecewise polynomial over time)'
$ws.Range("G5").Value = 'Cascaded PID (formation) + potential field (collision avoidance)'
$ws.Range("H5").Value = 'Agents'' state, obstacle pose'
$ws.Range("F6").Value = 'Payload trajectory set at beginning (polynomial over time), modified by dynamic motion primitives'
$ws.Range("G6").Value = 'Analytical feedback control strategy, dynamic motion primitives (collision avoidance)'
$ws.Range("H6").Value = 'Agents'' state, obstacle pose'
$ws.Range("F7").Value = 'Payload trajectory set at beginning'
$ws.Range("G7").Value = 'Admittance control'
$ws.Range("H7").Value = 'Agents'' state, cable force'
$ws.Range("F8").Value = 'Payload trajectory set at beginning'
$ws.Range("G8").Value = 'Force estimation + Admittance control + MPC'
$ws.Range("H8").Value = 'Agents'' state'
$ws.Range("B9").Value = 'Manipulation'
$ws.Range("F9").Value = 'Payload trajectory set at beginning'
$ws.Range("G9").Value = 'PID'
$ws.Range("H9").Value = 'Agents'' state'
$ws.Range("F10").Value = 'Payload trajectory set at beginning'
$ws.Range("G10").Value = 'Control force = PD + bias force'
$ws.Range("H10").Value = 'Internal force'
$ws.Range("D11").Value = 'Leader-follower structure'
$ws.Range("E11").Value = 'Distributed'
$ws.Range("F11").Value = 'Continual calculation of setpoint for group center position given desired payload pose and velocity'
$ws.Range("G11").Value = 'PD (trajectory) + LQR (formation) + PD (agent pose)'
$ws.Range("H11").Value = 'Agents'' state, load state'
$ws.Range("F12").Value = 'Payload trajectory set at beginning'
$ws.Range("G12").Value = 'Spring-damper (PD) Model or Potential field'
$ws.Range("H12").Value = 'Agents'' state'
$ws.Range("F13").Value = 'Payload path set at beginning, continual calculation of velocity'
$ws.Range("G13").Value = 'MPC & PI (velocity) + penalty force & VFH+- (collision avoidance)'
$ws.Range("H13").Value = 'Agents'' state with respect to formation, load pose, some communication with neighbors'
$ws.Range("F14").Value = 'Waypoints pre-specified, calculate reference trajectory subject to min. Snap (QP) & min. Duration (Coordinate descent)'
$ws.Range("G14").Value = 'Distributed wrench controller'
$ws.Range("H14").Value = 'Agents'' state, load state'
$ws.Range("F15").Value = 'Trajectory planning addressed in Lai, Wang, and Chen (2017) and Lai et al. (2016)'
$ws.Range("G15").Value = 'LQR (position) + Robust Perfect Tracking Control (trajectory)'
$ws.Range("H15").Value = 'Body state (Agents'' state)'
$ws.Range("F16").Value = 'Continual calculation of trajectory by leader given a desired payload final pose.'
$ws.Range("G16").Value = 'MPC receding horizon'
$ws.Range("H16").Value = 'Agents'' state, load geometry (no force/torque measurements)'
$ws.Range("F17").Value = 'Predefined trajectory tracking; Point-to-Point path planning'
$ws.Range("G17").Value = 'Robust optimal sliding mode control'
$ws.Range("H17").Value = 'Agents'' state'

$ws.Columns.Item(6).ColumnWidth = 22.85546875
